{"js": "// se escribio chao mundo parte2\n// Add a new paragraph with the text \"Chao mundo parte2\" right after the\n// existing \"Hola mundo parte1\" paragraph (i.e. at the end of the body).\n// Word automatically carries the surrounding run/paragraph formatting\n// (the es-ES language mark) onto the newly inserted paragraph.\nconst body = context.document.body;\nbody.insertParagraph(\"Chao mundo parte2\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Add a new paragraph \"Chao mundo parte2\" right after the existing\n# \"Hola mundo parte1\" paragraph, at the end of the document.\n$d = $word.ActiveDocument\n\n$endRange = $d.Content\n$endRange.Collapse(0)  # wdCollapseEnd\n$newPara = $endRange.InsertParagraphAfter()\n$endRange.Collapse(0)\n$endRange.InsertAfter(\"Chao mundo parte2\")\n"}
